# Reconciliation breaks combined - add custodian/organisation_name columns,
# refresh description/notes text, and update recommended_action + confidence.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new columns (custodian, organisation_name) right after
#    "bank_account" (column C), pushing classification..notes from D..J to F..L.
# ---------------------------------------------------------------------------
$ws.Columns("D:E").Insert()

# Explicit widths for the two new columns (match target layout).
$ws.Columns("D").ColumnWidth = 18.7109375
$ws.Columns("E").ColumnWidth = 18.7109375

# ---------------------------------------------------------------------------
# 2. Headers
# ---------------------------------------------------------------------------
$ws.Range("D1").Value = "custodian"
$ws.Range("E1").Value = "organisation_name"

# ---------------------------------------------------------------------------
# 3. New data for row 2 (priority 1 - Nestle SA / CHF quantity mismatch)
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = "CUST/UBSCH"
$ws.Range("E2").Value = "Nestle SA"

$ws.Range("I2").Value = "Custody cash reflects 12,000 shares (CHF 3.10 = CHF 37,200 gross; CHF 24,180 net) while NBIM booked 10,000 shares (CHF 31,000 gross; CHF 20,150 net). Dates, currency, and tax rate align; the mismatch is driven by entitlement quantity."

$ws.Range("J2").Value = 0.8

$ws.Range("K2").Value = "DRAFT_CUSTODIAN_TICKET"

$ws.Range("L2").Value = "Please confirm the entitled quantity at record date and whether 2,000 shares were ineligible (e.g., on loan or acquired ex-date). Custody record shows LOAN_QUANTITY 0 and LENDING_PERCENTAGE 0 but cash calculated on 12,000, while NOMINAL_BASIS = 10,000. Request breakdown of entitlement vs. holding and correct the paid amount or issue an adjustment if 10,000 is correct. If 12,000 is confirmed as entitled, NBIM will adjust its booking accordingly. POSSIBLE_RESTITUTION flags appear unrelated to this quantity issue."

# ---------------------------------------------------------------------------
# 4. New data for row 3 (priority 2 - Samsung Electronics / KRW tax mismatch)
# ---------------------------------------------------------------------------
$ws.Range("D3").Value = "CUST/HSBCKR"
$ws.Range("E3").Value = "Samsung Electronics Co Ltd"

$ws.Range("I3").Value = "NBIM applied ~25% total tax (22% WHT + ~3% local) while custody applied 20%, resulting in higher cash at custody. The USD difference (~342.77) aligns with ~5% of gross in KRW."

$ws.Range("J3").Value = 0.76

$ws.Range("K3").Value = "PROPOSE_NBIM_CORRECTION"

$ws.Range("L3").Value = "NBIM appears to have an extra local tax component (bringing total to ~25%). Please review KR dividend tax setup for this event and align to the rate applied by HSBC Korea (20%) or confirm the correct statutory/treaty rate. If 20% is confirmed, adjust NBIM booking and tax rates; if not, request custody to clarify their applied rate and any relief-at-source treatment."

Write-Output "Done."
